$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 181; this pushes the existing
# rows 181..236 down to 182..237 (matching the shift seen across the
# whole diff) and extends the used range to A1:R237.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new data record.
$ws.Cells.Item(181, 1).Value2 = 3
$ws.Cells.Item(181, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value2 = "Coquimbo"
$ws.Cells.Item(181, 4).Value2 = 44463
$ws.Cells.Item(181, 5).Value2 = 5
$ws.Cells.Item(181, 6).Value2 = 100112017
$ws.Cells.Item(181, 7).Value2 = "Apio"
$ws.Cells.Item(181, 8).Value2 = "Americana (o)"
$ws.Cells.Item(181, 9).Value2 = "Primera"
$ws.Cells.Item(181, 10).Value2 = 210
$ws.Cells.Item(181, 11).Value2 = 9000
$ws.Cells.Item(181, 12).Value2 = 9500
$ws.Cells.Item(181, 13).Value2 = 9262
$ws.Cells.Item(181, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(181, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(181, 16).Value2 = 1544
$ws.Cells.Item(181, 17).Value2 = 6
$ws.Cells.Item(181, 18).Value2 = "Hortaliza"

# Give the new row's date cell (column D) the same date-time number
# format used by every other row's D column.
$ws.Cells.Item(181, 4).NumberFormat = $ws.Cells.Item(182, 4).NumberFormat
